# formula: more math/trig functions
# Adds SIGN, SIN, SINH, SQRT, SQRTPI, SUMPRODUCT, SUMSQ, TAN, TANH, TRUNC
# test coverage to the "Math and Trig" worksheet of formulareference.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Math and Trig")

# ---- Header row (row 2): function name labels in new columns AW:BF ----
$ws.Range("AW2").Value = "SIGN"
$ws.Range("AX2").Value = "SIN"
$ws.Range("AY2").Value = "SINH"
$ws.Range("AZ2").Value = "SQRT"
$ws.Range("BA2").Value = "SQRTPI"
$ws.Range("BB2").Value = "SUMPRODUCT"
$ws.Range("BC2").Value = "SUMSQ"
$ws.Range("BD2").Value = "TAN"
$ws.Range("BE2").Value = "TANH"
$ws.Range("BF2").Value = "TRUNC"

# Match the header style (bold/centered) used by the rest of row 2
$ws.Range("AV2").Copy()
$ws.Range("AW2:BF2").PasteSpecial(-4122)  # xlPasteFormats

# ---- Helper data used by SUMPRODUCT / SUMSQ examples ----
$ws.Range("AX23").Value = 1
$ws.Range("AY23").Value = 2
$ws.Range("AX24").Value = 3
$ws.Range("AY24").Value = 4
$ws.Range("AX26").Value = 5
$ws.Range("AY26").Value = 7
$ws.Range("AX27").Value = 6
$ws.Range("AY27").Value = 8

# ---- Row 3 ----
$ws.Range("AW3").Formula = "=SIGN()"
$ws.Range("AX3").Formula = "=SIN(0)"
$ws.Range("AY3").Formula = "=LEFT(SINH(1),5)"
$ws.Range("AZ3").Formula = "=SQRT()"
$ws.Range("BA3").Formula = "=SQRTPI(81/PI())"
$ws.Range("BB3").Formula = "=SUMPRODUCT(AX23:AY24,AX26:AY27)"
$ws.Range("BC3").Formula = "=SUMSQ(4)"
$ws.Range("BD3").Formula = "=LEFT(TAN(45),6)"
$ws.Range("BE3").Formula = "=TANH(45)"
$ws.Range("BF3").Formula = "=TRUNC()"

# ---- Row 4 ----
$ws.Range("AW4").Formula = "=SIGN(0)"
$ws.Range("AX4").Formula = "=SIN(PI()/2)"
$ws.Range("AY4").Formula = "=SINH(0)"
$ws.Range("AZ4").Formula = "=SQRT(4)"
$ws.Range("BB4").Formula = "=SUMPRODUCT(1,2,3)"
$ws.Range("BC4").Formula = "=SUMSQ(AX23:AY23)"
$ws.Range("BE4").Formula = "=TANH(90)"
$ws.Range("BF4").Formula = "=TRUNC(8.9)"

# ---- Row 5 ----
$ws.Range("AW5").Formula = "=SIGN(-1.2)"
$ws.Range("AZ5").Formula = "=SQRT(81)"
$ws.Range("BB5").Formula = "=SUMPRODUCT(AX24:AY24,AX24:AY24)"
$ws.Range("BC5").Formula = "=SUMSQ(AX23:AY24)"
$ws.Range("BE5").Formula = "=LEFT(TANH(0.5),6)"
$ws.Range("BF5").Formula = "=TRUNC(-8.9)"

# ---- Row 6 ----
$ws.Range("AW6").Formula = "=SIGN(1.321)"
$ws.Range("AZ6").Formula = "=SQRT(-81)"
$ws.Range("BB6").Formula = "=SUMPRODUCT(AX24:AY24,AX26:AY26,AX27:AY27)"
$ws.Range("BF6").Formula = "=TRUNC(0.45)"

# ---- Row 7 ----
$ws.Range("AW7").Formula = '=SIGN("A")'
$ws.Range("BF7").Formula = "=TRUNC(1.23,2)"

# ---- Row 8 ----
$ws.Range("BF8").Formula = "=TRUNC(1.23,1)"

# ---- Row 9 ----
$ws.Range("BF9").Formula = "=TRUNC(-1.23,1)"

# ---- Row 10 ----
$ws.Range("BF10").Formula = "=TRUNC(-1.23,4)"

# ---- Row 11 ----
$ws.Range("BF11").Formula = "=TRUNC(1.23,-2)"

# ---- Row 12 ----
$ws.Range("BF12").Formula = "=TRUNC(1.23,-5)"

# ---- Selection / scroll position, matching the commit's final cursor spot ----
$ws.Range("BD23").Select()
